$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-29 10:50:54"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
